$wb = $excel.ActiveWorkbook

$siteGroup = $wb.Worksheets.Item("Site Group")

# The "region_id" (col B) and "site_group_name" (col C) columns were swapped
# for convention. Use a native cut/insert column move so values, shared
# strings and cell styles all travel together exactly like a manual
# drag-and-drop column reorder in Excel.
$colB = $siteGroup.Columns.Item(2)
$colC = $siteGroup.Columns.Item(3)
$colC.Cut() | Out-Null
$colB.Insert() | Out-Null

# The "Site Group" sheet becomes the active sheet/tab, with a new selection.
$siteGroup.Range("E5").Select() | Out-Null
$siteGroup.Activate() | Out-Null
